$wb = $excel.ActiveWorkbook

# Update status text "Ready for handoff" -> "In Translation" on all three sheets
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow columns E:F on Overview sheet and column C on zh-cn / de-de sheets
# (target OOXML width = 13.4101845877511 "characters"; the COM layer here only
# persists ColumnWidth on a 1/6-character grid, so 12.5 is the input that lands
# on the closest achievable grid point, 13.333333333333334)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
